$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'43.207.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "'2.323.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'303.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'102.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.24%  "
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.92%  "
$ws.Range("D10").Value = "'36.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.65%  "
$ws.Range("D11").Value = "'0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("D13").Value = "'17.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.03%  "
$ws.Range("D14").Value = "'6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "'2.674.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "'2.320.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "'0.811"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "'43.097.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.87%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.16%  "
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "'67.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'237.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +13.80%  "
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'24.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.75%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'34.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'168.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "'4.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").Value = "'5.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").Value = "'17.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.21%  "
$ws.Range("E36").Value = "  +3.27%  "
$ws.Range("D37").Value = "'0.0695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").Value = "'2.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  +4.92%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.993.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("D44").Value = "'0.0290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("D45").Value = "'10.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.22%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.01%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'17.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "'56.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.09%  "
$ws.Range("D49").Value = "'2.546.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("E50").Value = "  +3.65%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.29%  "
